# Commit: "added selection of xls sheet_name, fixed PyPi long Description"
#
# Functional change: add a second worksheet ("2nd Sheet") after the existing
# "DataDriven" sheet. It reuses the same *** Test Cases *** / ${username} /
# ${password} / [Tags] / [Documentation] header row, then rows with the
# numbers 1..8 (stored as text, matching column B) in column B and the
# literal "Test" in column C - illustrating DataDriver's ability to pick a
# named Excel sheet ("sheet_name") instead of just the first one.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- add the new sheet right after "DataDriven" -----------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2nd Sheet"

# Store everything as text (same presentation as sheet1, which applies a
# "@" text number format to its whole used range).
$ws2.Range("A1:E9").NumberFormat = "@"
$ws2.Range("B10:C10").NumberFormat = "@"

# --- header row, identical wording to "DataDriven" ---------------------
$ws2.Range("A1").Value = "*** Test Cases ***"
$ws2.Range("B1").Value = "`${username}"
$ws2.Range("C1").Value = "`${password}"
$ws2.Range("D1").Value = "[Tags]"
$ws2.Range("E1").Value = "[Documentation]"

# --- data rows: B2:B9 = 1..8 (text), C2:C9 = "Test" ---------------------
for ($i = 1; $i -le 8; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 2).Value = [string]$i
}
for ($i = 1; $i -le 8; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 3).Value = "Test"
}

# Match column widths to the content (best effort / cosmetic).
$ws2.Range("A1:E10").EntireColumn.AutoFit()

# Same page setup as "DataDriven" (A4/Letter #9, portrait).
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# --- selection / active-sheet bookkeeping -------------------------------
# Previously sheet1 was the tab-selected / active sheet with a stray
# selection left at D16. Select the whole used range on sheet1 instead,
# then make the new "2nd Sheet" the active tab.
$ws1.Activate() | Out-Null
$ws1.Range("A1:E9").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B20").Select() | Out-Null
